$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, pushing existing rows 9..121 down to 10..122.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new data record
# (same market/product metadata as its neighbours, new variety/price details).
$ws.Range("A9").Value = 10
$ws.Range("B9").Value = "Vega Modelo de Temuco"
$ws.Range("C9").Value = "La Araucanía"
$ws.Range("D9").Value = 44515
$ws.Range("E9").Value = 9
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100103
$ws.Range("H9").Value = "Frutos de hueso (carozo)"
$ws.Range("I9").Value = 100103001
$ws.Range("J9").Value = "Cereza"
$ws.Range("K9").Value = "Early Burlat"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 30
$ws.Range("N9").Value = 3500
$ws.Range("O9").Value = 3500
$ws.Range("P9").Value = 3500
$ws.Range("Q9").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R9").Value = "Región del Maule"
$ws.Range("S9").Value = 3500
$ws.Range("T9").Value = 1

# Match the date-formatted style used by the other rows' Fecha column (D).
$ws.Range("D9").NumberFormat = $ws.Range("D10").NumberFormat()
